$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the greeting text for the R10 rule (row 8, column E)
$ws.Range("E8").Value = "GIT UPDATE"

# Match the author's cursor position / selection recorded in the file
$ws.Range("E8").Select()
